$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 29288.428
$ws.Range("I21").Value = 31254.75
$ws.Range("J21").Value = 26666.666
$ws.Range("K21").Value = 31254.75
$ws.Range("L21").Value = 26666.666
$ws.Range("M21").Value = -30786.75
$ws.Range("N21").Value = -27602.666

$ws.Range("H23").Value = 29288.428
$ws.Range("I23").Value = 31254.75
$ws.Range("J23").Value = 26666.666
$ws.Range("K23").Value = 31254.75
$ws.Range("L23").Value = 26666.666
$ws.Range("M23").Value = -31020.75
$ws.Range("N23").Value = -27134.666

$ws.Range("H29").Value = 1134.3334
$ws.Range("J29").Value = 3000
$ws.Range("L29").Value = 9000
$ws.Range("N29").Value = -9562

$ws.Range("H75").Value = 27465.273
$ws.Range("J75").Value = 27465.273
$ws.Range("L75").Value = 27465.273
$ws.Range("N75").Value = -29337.273

$ws.Range("H78").Value = 27465.273
$ws.Range("J78").Value = 27465.273
$ws.Range("L78").Value = 82395.819
$ws.Range("N78").Value = -91755.819

$ws.Range("H138").Value = 2803.88
$ws.Range("I138").Value = 1145.762
$ws.Range("J138").Value = 3244.6455
$ws.Range("K138").Value = 3437.286
$ws.Range("L138").Value = 9733.9365
$ws.Range("M138").Value = 1702.714
$ws.Range("N138").Value = -20013.9365

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 17333.334
$ws.Range("J9").Value = 17333.334
$ws.Range("L9").Value = 17333.334
$ws.Range("N9").Value = -17673.334

$ws.Range("H20").Value = 17333.334
$ws.Range("J20").Value = 17333.334
$ws.Range("L20").Value = 17333.334
$ws.Range("N20").Value = -17873.334

$ws.Range("H23").Value = 10927.714
$ws.Range("I23").Value = 11248.5
$ws.Range("J23").Value = 10500
$ws.Range("K23").Value = 11248.5
$ws.Range("L23").Value = 10500
$ws.Range("M23").Value = -10989.5
$ws.Range("N23").Value = -11018

$ws.Range("H32").Value = 14946583
$ws.Range("I32").Value = 25010844
$ws.Range("J32").Value = 36568.15
$ws.Range("K32").Value = 25010844
$ws.Range("L32").Value = 36568.15
$ws.Range("M32").Value = -25010557
$ws.Range("N32").Value = -37142.15

$ws.Range("H50").Value = 520.4286
$ws.Range("I50").Value = 196.5
$ws.Range("J50").Value = 650
$ws.Range("K50").Value = 196.5
$ws.Range("L50").Value = 650
$ws.Range("M50").Value = 517.5
$ws.Range("N50").Value = -2078

$ws.Range("H54").Value = 19999.5
$ws.Range("J54").Value = 19999.5
$ws.Range("L54").Value = 19999.5
$ws.Range("N54").Value = -21537.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 26600
$ws.Range("J15").Value = 26600
$ws.Range("L15").Value = 26600
$ws.Range("N15").Value = -27054

$ws.Range("H81").Value = 7793.3335
$ws.Range("J81").Value = 7793.3335
$ws.Range("L81").Value = 7793.3335
$ws.Range("N81").Value = -9915.333500000001

$ws.Range("H84").Value = 7793.3335
$ws.Range("J84").Value = 7793.3335
$ws.Range("L84").Value = 23380.0005
$ws.Range("N84").Value = -33988.00049999999

$ws.Range("H135").Value = 29700
$ws.Range("J135").Value = 29700
$ws.Range("L135").Value = 29700
$ws.Range("N135").Value = -39840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H102").Value = 27300
$ws.Range("J102").Value = 27300
$ws.Range("L102").Value = 27300
$ws.Range("N102").Value = -32168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 602.0816
$ws.Range("I5").Value = 488.14285
$ws.Range("K5").Value = 1464.42855
$ws.Range("M5").Value = -1352.42855

$ws.Range("H87").Value = 22655.24
$ws.Range("I87").Value = 5866
$ws.Range("J87").Value = 30556.059
$ws.Range("K87").Value = 17598
$ws.Range("L87").Value = 91668.177
$ws.Range("M87").Value = -16350
$ws.Range("N87").Value = -94164.177

$ws.Range("H90").Value = 22655.24
$ws.Range("I90").Value = 5866
$ws.Range("J90").Value = 30556.059
$ws.Range("K90").Value = 52794
$ws.Range("L90").Value = 275004.531
$ws.Range("M90").Value = -46554
$ws.Range("N90").Value = -287484.531

$ws.Range("H121").Value = 1842.0476
$ws.Range("I121").Value = 1300
$ws.Range("J121").Value = 1899.1052
$ws.Range("K121").Value = 3900
$ws.Range("L121").Value = 5697.3156
$ws.Range("M121").Value = -2590
$ws.Range("N121").Value = -8317.3156

$ws.Range("H131").Value = 744
$ws.Range("I131").Value = 307.27274
$ws.Range("J131").Value = 909.65515
$ws.Range("K131").Value = 921.81822
$ws.Range("L131").Value = 2728.96545
$ws.Range("M131").Value = 4118.18178
$ws.Range("N131").Value = -12808.96545

$ws.Range("H135").Value = 602.0816
$ws.Range("I135").Value = 488.14285
$ws.Range("K135").Value = 4393.28565
$ws.Range("M135").Value = -1858.28565

$ws.Range("H137").Value = 4247.0586
$ws.Range("I137").Value = 3390.9092
$ws.Range("J137").Value = 5816.6665
$ws.Range("K137").Value = 10172.7276
$ws.Range("L137").Value = 17449.9995
$ws.Range("M137").Value = -5072.7276
$ws.Range("N137").Value = -27649.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1421.2
$ws.Range("I102").Value = 1270.25
$ws.Range("J102").Value = 2025
$ws.Range("K102").Value = 1270.25
$ws.Range("L102").Value = 2025
$ws.Range("M102").Value = 351.75
$ws.Range("N102").Value = -5269

$ws.Range("H122").Value = 5560604.5
$ws.Range("I122").Value = 14287229
$ws.Range("J122").Value = 7298
$ws.Range("K122").Value = 42861687
$ws.Range("L122").Value = 21894
$ws.Range("M122").Value = -42859237
$ws.Range("N122").Value = -26794

$ws.Range("H133").Value = 43386.332
$ws.Range("J133").Value = 43386.332
$ws.Range("L133").Value = 43386.332
$ws.Range("N133").Value = -53506.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4039.9473
$ws.Range("I40").Value = 3480.6924
$ws.Range("J40").Value = 5251.6665
$ws.Range("K40").Value = 3480.6924
$ws.Range("L40").Value = 5251.6665
$ws.Range("M40").Value = -3344.6924
$ws.Range("N40").Value = -5523.6665

$ws.Range("H51").Value = 19389.334
$ws.Range("J51").Value = 19389.334
$ws.Range("L51").Value = 19389.334
$ws.Range("N51").Value = -20345.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4178.564
$ws.Range("I136").Value = 4498.393
$ws.Range("J136").Value = 3364.4546
$ws.Range("K136").Value = 13495.179
$ws.Range("L136").Value = 10093.3638
$ws.Range("M136").Value = -10945.179
$ws.Range("N136").Value = -15193.3638
